# Stock hasta el punto pedido por la empresa
# Replace the placeholder "item_N" codes with the real MCER item codes,
# re-shuffle the stock quantities on packing_1, and trim both sheets
# down to 27 data rows (A1:B28) instead of 30 (A1:B31).

$wb = $excel.ActiveWorkbook

$items = @(
    "MCER017","MCER018","MCER020","MCER021","MCER022","MCER026","MCER027",
    "MCER028","MCER029","MCER030","MCER031","MCER032","MCER033","MCER034",
    "MCER035","MCER036","MCER037","MCER038","MCER039","MCER040","MCER041",
    "MCER043","MCER046","MCER047","MCER051","MCER057","MCER067"
)

$stock0 = @(918,291,202,797,117,551,115,343,773,894,879,563,300,695,940,656,292,704,348,495,64,561,565,147,158,43,292)
$stock1 = @(890,837,821,177,988,499,758,939,602,250,170,188,890,149,542,223,889,800,146,805,973,648,850,818,710,439,16)

$ws0 = $wb.Worksheets.Item("packing_0")
$ws1 = $wb.Worksheets.Item("packing_1")

for ($i = 0; $i -lt $items.Length; $i++) {
    $row = $i + 2
    $ws0.Cells.Item($row, 1).Value = $items[$i]
    $ws0.Cells.Item($row, 2).Value = $stock0[$i]
    $ws1.Cells.Item($row, 1).Value = $items[$i]
    $ws1.Cells.Item($row, 2).Value = $stock1[$i]
}

# Drop the now-unused trailing rows (29, 30, 31) from both sheets so the
# used range shrinks back down to A1:B28.
$ws0.Rows.Item(31).Delete() | Out-Null
$ws0.Rows.Item(30).Delete() | Out-Null
$ws0.Rows.Item(29).Delete() | Out-Null

$ws1.Rows.Item(31).Delete() | Out-Null
$ws1.Rows.Item(30).Delete() | Out-Null
$ws1.Rows.Item(29).Delete() | Out-Null
